$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended at the bottom of the table (dates were being
# mis-parsed before; 44217 = 2021-01-21, 44220 = 2021-01-24).
$ws.Range("A16").Value = 44217
$ws.Range("B16").Value = 37960000
$ws.Range("C16").Value = "'"
$ws.Range("D16").Value = "'"
$ws.Range("E16").Value = 17546374
$ws.Range("F16").Value = 15053257
$ws.Range("G16").Value = 2394961
$ws.Range("H16").Value = 2089181

$ws.Range("A17").Value = 44220
$ws.Range("B17").Value = 41411550
$ws.Range("C17").Value = "'"
$ws.Range("D17").Value = "'"
$ws.Range("E17").Value = 21848655
$ws.Range("F17").Value = 18502131
$ws.Range("G17").Value = 3216836
$ws.Range("H17").Value = 2567018

# Match the date-format styling used by the rest of column A.
$ws.Range("A16").NumberFormat = $ws.Range("A15").NumberFormat
$ws.Range("A17").NumberFormat = $ws.Range("A15").NumberFormat

# The "'" placeholder above (mirroring the blank C/D cells already present
# in rows 12-15) leaves a stray quote-prefix format on the cell; clear it
# by pasting the plain number format from column B over C:D so the new
# blank cells carry the same unstyled formatting as the existing ones.
$ws.Range("B16").Copy() | Out-Null
$ws.Range("C16:D16").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Copy() | Out-Null
$ws.Range("C17:D17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
